$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-typed assignment (leading apostrophe = Excel quote-prefix)
# so numeric-looking strings like "1.010" or "29.705.93" keep their exact
# literal text instead of being parsed into numbers.
$ws.Range('D2').Value = "'29.705.93"
$ws.Range('E2').Value = "'  -2.87%  "
$ws.Range('D3').Value = "'2.088.59"
$ws.Range('E3').Value = "'  -1.17%  "
$ws.Range('D4').Value = "'1.010"
$ws.Range('E4').Value = "'  -0.30%  "
$ws.Range('D5').Value = "'345.26"
$ws.Range('E5').Value = "'  -0.28%  "
$ws.Range('D6').Value = "'1.009"
$ws.Range('E6').Value = "'  -0.20%  "
$ws.Range('D7').Value = "'0.5166"
$ws.Range('E7').Value = "'  -1.73%  "
$ws.Range('D8').Value = "'0.4400"
$ws.Range('E8').Value = "'  -2.56%  "
$ws.Range('D9').Value = "'0.09296"
$ws.Range('E9').Value = "'  +3.16%  "
$ws.Range('D10').Value = "'51.58"
$ws.Range('E10').Value = "'  -4.48%  "
$ws.Range('D11').Value = "'1.178"
$ws.Range('E11').Value = "'  +0.63%  "
$ws.Range('D12').Value = "'25.44"
$ws.Range('E12').Value = "'  +4.34%  "
$ws.Range('D13').Value = "'2.095.08"
$ws.Range('E13').Value = "'  -0.76%  "
$ws.Range('D14').Value = "'6.746"
$ws.Range('E14').Value = "'  -0.85%  "
$ws.Range('D15').Value = "'8.132"
$ws.Range('E15').Value = "'  +0.73%  "
$ws.Range('D16').Value = "'99.52"
$ws.Range('E16').Value = "'  -0.17%  "
$ws.Range('D17').Value = "'0.00001166"
$ws.Range('E17').Value = "'  -0.80%  "
$ws.Range('D18').Value = "'1.010"
$ws.Range('E18').Value = "'  -0.32%  "
$ws.Range('E19').Value = "'  +8.46%  "
$ws.Range('D20').Value = "'0.06674"
$ws.Range('E20').Value = "'  -0.48%  "
$ws.Range('D21').Value = "'1.006"
$ws.Range('E21').Value = "'  -0.56%  "
$ws.Range('D22').Value = "'6.194"
$ws.Range('E22').Value = "'  -2.03%  "
$ws.Range('D23').Value = "'29.781.09"
$ws.Range('E23').Value = "'  -2.87%  "
$ws.Range('D24').Value = "'12.72"
$ws.Range('E24').Value = "'  -0.48%  "
$ws.Range('D25').Value = "'2.307"
$ws.Range('E25').Value = "'  -3.56%  "
$ws.Range('D26').Value = "'2.339.75"
$ws.Range('E26').Value = "'  -0.89%  "
$ws.Range('D27').Value = "'21.92"
$ws.Range('E27').Value = "'  -1.55%  "
$ws.Range('D28').Value = "'163.27"
$ws.Range('E28').Value = "'  -1.19%  "
$ws.Range('D29').Value = "'2.528"
$ws.Range('E29').Value = "'  -0.09%  "
$ws.Range('D30').Value = "'132.87"
$ws.Range('E30').Value = "'  -1.57%  "
$ws.Range('D31').Value = "'1.147"
$ws.Range('E31').Value = "'  -3.87%  "
$ws.Range('D32').Value = "'0.1054"
$ws.Range('E32').Value = "'  -1.59%  "
$ws.Range('D33').Value = "'1.619"
$ws.Range('E33').Value = "'  -0.99%  "
$ws.Range('D34').Value = "'6.203"
$ws.Range('E34').Value = "'  -2.44%  "
$ws.Range('D35').Value = "'3.953"
$ws.Range('E35').Value = "'  -0.74%  "
$ws.Range('D36').Value = "'6.167"
$ws.Range('E36').Value = "'  +4.52%  "
$ws.Range('D37').Value = "'10.14"
$ws.Range('E37').Value = "'  -0.56%  "
$ws.Range('D38').Value = "'0.02579"
$ws.Range('E38').Value = "'  -2.27%  "
$ws.Range('D39').Value = "'0.06739"
$ws.Range('E39').Value = "'  -1.44%  "
$ws.Range('D40').Value = "'0.2281"
$ws.Range('E40').Value = "'  -1.95%  "
$ws.Range('D41').Value = "'0.6871"
$ws.Range('E41').Value = "'  +0.11%  "
$ws.Range('D42').Value = "'12.39"
$ws.Range('E42').Value = "'  -1.98%  "
$ws.Range('D43').Value = "'1.284"
$ws.Range('E43').Value = "'  +1.36%  "
$ws.Range('D44').Value = "'0.6653"
$ws.Range('E44').Value = "'  +3.59%  "
$ws.Range('D45').Value = "'14.09"
$ws.Range('E45').Value = "'  -5.61%  "
$ws.Range('E46').Value = "'  -0.24%  "
$ws.Range('D47').Value = "'3.624"
$ws.Range('E47').Value = "'  -3.20%  "
$ws.Range('D48').Value = "'1.218"
$ws.Range('E48').Value = "'  -2.86%  "
$ws.Range('D49').Value = "'0.00000000340"
$ws.Range('E49').Value = "'  -6.50%  "
$ws.Range('D50').Value = "'81.91"
$ws.Range('D51').Value = "'0.07102"
$ws.Range('E51').Value = "'  -2.48%  "
